# Update countries & provincias Spain
# Refresh the COVID-19 snapshot: new timestamp, new case/death/etc. figures
# for a set of countries, and re-rank three pairs of countries whose case
# counts crossed over (Francia/Colombia, Irlanda/Kirguistan,
# Mozambique/Maldivas) by swapping their row labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 20:41"

# Country-label swaps caused by re-sorting on total cases (column B)
$ws.Range("A10").Value = "Francia"
$ws.Range("A11").Value = "Colombia"

$ws.Range("A68").Value = "Irlanda"
$ws.Range("A69").Value = "Kirguistan"

$ws.Range("A107").Value = "Mozambique"
$ws.Range("A108").Value = "Maldivas"

# Updated metrics: row number -> column letter -> new value
$changes = @{
    4   = @{ B = 8618942;  C = 34123;  D = 5617171; E = 2773987; G = 376; H = 227784 }
    5   = @{ B = 7758713;  C = 53555;  D = 6944258; E = 697125;  G = 677; H = 117330 }
    10  = @{ B = 999043;   C = 41622;  D = 108599;  E = 856234;  G = 162; H = 34210 }
    11  = @{ B = 981700;               D = 884895;  E = 67341;            H = 29464 }
    20  = @{ B = 399315;   C = 7960;                E = 87183;   G = 33;  H = 10032 }
    35  = @{ B = 186731;   C = 4151;   D = 154481;  E = 29118;   G = 53;  H = 3132 }
    62  = @{ B = 67027;    C = 1450;   D = 31409;   E = 35066;   G = 16;  H = 552 }
    67  = @{ B = 55357;    C = 276;    D = 38618;   E = 14851;   G = 8;   H = 1888 }
    68  = @{ B = 54476;    C = 1054;   D = 23364;   E = 29241;   G = 3;   H = 1871 }
    69  = @{ B = 54006;    C = 547;    D = 46726;   E = 6158;    G = 4;   H = 1122 }
    104 = @{ B = 12460;    C = 54;     D = 10609;   E = 1718 }
    107 = @{ B = 11559;    C = 228;    D = 9226;    E = 2252;    G = 2;   H = 81 }
    108 = @{ B = 11358;    C = 42;     D = 10383;   E = 938;              H = 37 }
    112 = @{ B = 10342;    C = 47;                  E = 278 }
    125 = @{ B = 5874;     C = 10;     D = 4764;    E = 927 }
    151 = @{ B = 3440;     C = 12;     D = 2608;    E = 700 }
    152 = @{ B = 3154;     C = 188;                 E = 1685 }
    157 = @{ B = 2414;     C = 8;      D = 1869;    E = 480 }
    166 = @{                          D = 1278;    E = 25 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
